# NV-37 Bác Sĩ Thảo 7-2024 — "Đơn 1 bác sĩ" sheet restructuring.
# The table is reshaped from the full (sale-chinh/sale-phu + doctor/phu-phau
# breakdown) layout out to a simplified "1 bác sĩ" layout: columns G..N are
# redefined and columns O..AA (Trả sau, Dư nợ, Bác sĩ 1/2, Phụ phẫu 1/2,
# Công phụ phẫu 1/2, Tỉ lệ/Chiết khấu sale chính/phụ, ...) are dropped in
# favour of two new "bác sĩ 1" discount-rate / discount-amount columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn 1 bác sĩ")

# --- Drop the now-unused trailing columns (O:AA) for rows 1-4 -------------
# Clearing (rather than deleting cells/shifting) lets the used range shrink
# naturally down to A1:N4, matching the updated <dimension>.
$ws.Range("O1:AA4").ClearContents()

# --- Header row (row 1): relabel G:N ---------------------------------------
$ws.Range("G1").Value = "Tên dịch vụ"
$ws.Range("H1").Value = "Đơn giá gốc"
$ws.Range("I1").Value = "Sale phụ"
$ws.Range("J1").Value = "Upsale"
$ws.Range("K1").Value = "Đơn giá"
$ws.Range("L1").Value = "Đã thanh toán"
$ws.Range("M1").Value = "Tỉ lệ chiết khấu bác sĩ 1"
$ws.Range("N1").Value = "Chiết khấu bác sĩ 1"

# --- Row 2 data (service: Treo cung) ---------------------------------------
$ws.Range("G2").Value = "Treo cung"
$ws.Range("H2").Value = 25000000
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = 25000000
$ws.Range("L2").Value = 25000000
$ws.Range("M2").Value = 0.1
$ws.Range("N2").Value = 2500000

# --- Row 3 data (service: cắt sẹo) ------------------------------------------
$ws.Range("G3").Value = "cắt sẹo "
$ws.Range("H3").Value = 7000000
$ws.Range("I3").Value = "Lê Đình Hậu"
$ws.Range("J3").Value = 6000000
$ws.Range("K3").Value = 13000000
$ws.Range("L3").Value = 13000000
$ws.Range("M3").Value = 0.1
$ws.Range("N3").Value = 1300000

# --- Row 4 data (totals row) ------------------------------------------------
# G4/I4 are already blank in the source layout (they were blank under the
# old column mapping too) — leave them untouched rather than re-writing an
# empty string, so they keep their existing empty-cell representation.
$ws.Range("H4").Value = 32000000
$ws.Range("J4").Value = 6000000
$ws.Range("K4").Value = 38000000
$ws.Range("L4").Value = 38000000
$ws.Range("M4").Value = 0.2
$ws.Range("N4").Value = 3800000
